$d = $word.ActiveDocument

# --- Change 1: FirstParagraph - NOAA wording update ---
$d.Content.Find.Execute(
    "National Oceanic and Atmospheric Administration (NOAA) Global Historical Climatology Network",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "NOAA National Centers for Environmental Information (NCEI) Global Historical Climatology Network",
    2) | Out-Null

# --- Change 2: split the "Including an interaction..." paragraph and insert
#     the new Vancouver/AccuWeather paragraph between the two halves ---
$vancouverPara = "A complication arose for Vancouver, as the NOAA NCEI station did not report any weather data from September 2025 onward. Since the model relies on October-to-present climate data, this posed a challenge. To address this, publicly available AccuWeather data provided by the course instructor was used instead. Because the AccuWeather records only included January" + [char]0x2013 + "February 2026, the model for Vancouver was adjusted to use this limited window. While precipitation data was not available in this substitute source, the impact is expected to be minimal given that precipitation was not a significant predictor in the larger multi-location model. This approach still allowed tracking of temperature conditions relevant to bloom timing, particularly exposure to extreme cold."

$replacement2 = "temporal shifts.^p" + $vancouverPara + "^pResults indicate"

$d.Content.Find.Execute(
    "temporal shifts. Results indicate",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $replacement2,
    2) | Out-Null

# --- Change 3: split the final paragraph, dropping the stray single-space
#     run that used to sit between the two sentences ---
$d.Content.Find.Execute(
    "for submission. Overall, this analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "for submission.^pOverall, this analysis",
    2) | Out-Null
